$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) "Programming" hours: D8 changes from 6 to 9
# ---------------------------------------------------------------------------
$ws.Range("D8").Value = 9

# ---------------------------------------------------------------------------
# 2) The "Content definieren" block (E9:E11 / F9:F11, each merged across 3
#    rows) shrinks to a 2-row block (E9:E10 / F9:F10), freeing up row 11
#    for a new "Git" entry with the value 3.
# ---------------------------------------------------------------------------

# Unmerge the old 3-row blocks so each row becomes independently addressable.
$ws.Range("E9:E11").UnMerge()
$ws.Range("F9:F11").UnMerge()

# Re-merge the top two rows into the new (shorter) block.
$ws.Range("E9:E10").Merge()
$ws.Range("F9:F10").Merge()

# The border formatting used to "close" a merged block moves down from row 11
# to row 10, and row 11 reverts to the plain "middle" formatting (column E is
# filled, column F is not). Swap the two rows' formats using a scratch cell.
$ws.Range("E11").Copy()
$ws.Range("Z1").PasteSpecial(-4122)

$ws.Range("E10").Copy()
$ws.Range("E11").PasteSpecial(-4122)

$ws.Range("Z1").Copy()
$ws.Range("E10").PasteSpecial(-4122)

$ws.Range("F11").Copy()
$ws.Range("F10").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false
$ws.Range("Z1").Clear()

# New content for the freed-up row.
$ws.Range("E11").Value = "Git"
$ws.Range("F11").Value = 3

# ---------------------------------------------------------------------------
# 3) Selection state, as last saved by the author.
# ---------------------------------------------------------------------------
$ws.Range("F1:F13").Select()
